$d = $word.ActiveDocument

# Find the paragraph ending with "... should also be here)." and insert
# a new list item right after it (before the "Source and documents" heading).
$target = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text -match "should also be here\)\.") {
        $target = $p
        break
    }
}

if ($target -eq $null) {
    throw "Could not find target paragraph"
}

# Insert a brand new paragraph right after this one; it inherits the same
# list/paragraph formatting (Listeafsnit style + numbering).
$target.Range.InsertParagraphAfter()

# The newly created paragraph is the one after $target now.
$newPara = $target.Next()
$newPara.Range.Text = "Add a few more queries"
